$d = $word.ActiveDocument

# Locate the run of text "For example " that needs to be split into three
# runs: "For " / "example," / " " (i.e. a comma is inserted after
# "example").
$found = $d.Content
$found.Find.Execute("For example ", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)

$start = $found.Start

# Insert the comma right after "example" (offset 11 within the match,
# i.e. right before the trailing space).
$commaPoint = $d.Range($start + 11, $start + 11)
$commaPoint.InsertAfter(",")

# Text is now "For example, " - force a run boundary between "For " and
# "example," by toggling a formatting property on the latter span and
# reverting it (same net formatting, but a new run is materialized).
$exampleComma = $d.Range($start + 4, $start + 12)
$exampleComma.Font.Bold = $true
$exampleComma.Font.Bold = $false

# Force a run boundary between "example," and the trailing space the same
# way.
$trailingSpace = $d.Range($start + 12, $start + 13)
$trailingSpace.Font.Bold = $true
$trailingSpace.Font.Bold = $false
